$d = $word.ActiveDocument

# --- Step 1: locate the paragraph that starts the "Usability" bullet list
# ("Record tutorial videos...") - the 6 new bullet points are inserted
# directly before it.
$targetIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Record tutorial videos")) {
        $targetIndex = $i
        break
    }
    $i = $i + 1
}
if ($targetIndex -eq -1) {
    throw "Could not find target paragraph 'Record tutorial videos...'"
}

$target = $d.Paragraphs($targetIndex)
$insertPos = $target.Range.Start
$insertRange = $d.Range($insertPos, $insertPos)

$newParasXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t xml:space="preserve">Add question to survey that asks whether patient has been diagnosed with cancer. If yes, do survey. </w:t></w:r><w:r><w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>Otherwise,</w:t></w:r><w:r><w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t xml:space="preserve"> don’t</w:t></w:r><w:r><w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>Individual responses to survey as well as score to send to clinician.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>Add delete button to photos – but warn them against it</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>Default is email far shot, first mole, last mole. User can then add risk assessment and other photos.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>Worth looking about AI models to classify tumours</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:lastRenderedPageBreak/><w:t>Patient mental health patients with bad outome may be suicidal. Don’t make app tell patients its bad. Leave that to clinician.</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $insertRange.InsertXML($newParasXml)

# InsertXML merges the formatting of its *last* paragraph mark into the
# paragraph that followed the insertion point, which is why the fragment
# above carries one extra trailing empty <w:p/> to absorb that merge.
# That now-empty paragraph sits right before the original target
# paragraph ("Record tutorial videos...") and must be removed.
$artifactIndex = $targetIndex + 6
$artifact = $d.Paragraphs($artifactIndex)
if ($artifact.Range.Text -ne "") {
    throw "Unexpected artifact paragraph content: $($artifact.Range.Text)"
}
$artifact.Range.Delete()

# --- Step 2: remove the <w:lastRenderedPageBreak/> marker from the
# "Email feature not working..." paragraph (it moved earlier in the
# document because of the newly inserted page-break marker above).
$emailIndex = -1
$i = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Email feature not working")) {
        $emailIndex = $i
        break
    }
    $i = $i + 1
}
if ($emailIndex -eq -1) {
    throw "Could not find target paragraph 'Email feature not working...'"
}

$emailPara = $d.Paragraphs($emailIndex)
$emailRange = $emailPara.Range
$replacementXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="6191AEF4" w14:textId="779B5DEA" w:rsidR="004C4386" w:rsidRDefault="003A3949" w:rsidP="003A3949"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/></w:rPr></w:pPr><w:r w:rsidRPr="003A3949"><w:rPr><w:color w:val="538135" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>Email feature not working in app currently. This needs to be fixed. Last year’s team said they found this very difficult.</w:t></w:r></w:p>'
$null = $emailRange.InsertXML($replacementXml)
